$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain plain text even when the value looks numeric,
# by applying a text number format before assigning the value. This avoids Excel
# auto-converting strings like "408.40" or "0.0000214" into floating point numbers.
$priceCells = @(
    "D2", "D3", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D36", "D37", "D39", "D40", "D43", "D46", "D47", "D48", "D49", "D51"
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "61.931.97"
$ws.Range("E2").Value = "  +0.65%  "

# Row 3
$ws.Range("D3").Value = "3.432.41"
$ws.Range("E3").Value = "  +1.60%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "408.40"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6
$ws.Range("D6").Value = "128.44"
$ws.Range("E6").Value = "  -3.71%  "

# Row 7
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +5.88%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  +9.83%  "

# Row 10
$ws.Range("E10").Value = "  +20.46%  "

# Row 11
$ws.Range("D11").Value = "42.44"
$ws.Range("E11").Value = "  +0.17%  "

# Row 12
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").Value = "0.0000214"
$ws.Range("E12").Value = "  +67.89%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.141"
$ws.Range("E13").Value = "  -0.32%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.969.41"
$ws.Range("E14").Value = "  +1.45%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "21.31"
$ws.Range("E15").Value = "  +7.96%  "

# Row 16
$ws.Range("D16").Value = "8.90"
$ws.Range("E16").Value = "  +5.95%  "

# Row 17
$ws.Range("D17").Value = "3.415.09"
$ws.Range("E17").Value = "  +1.20%  "

# Row 18
$ws.Range("D18").Value = "12.35"
$ws.Range("E18").Value = "  +12.79%  "

# Row 19
$ws.Range("D19").Value = "1.08"
$ws.Range("E19").Value = "  +5.98%  "

# Row 20
$ws.Range("D20").Value = "61.949.39"
$ws.Range("E20").Value = "  +0.61%  "

# Row 21
$ws.Range("D21").Value = "398.89"
$ws.Range("E21").Value = "  +25.91%  "

# Row 22
$ws.Range("E22").Value = "  +5.95%  "

# Row 23
$ws.Range("D23").Value = "3.20"
$ws.Range("E23").Value = "  +0.39%  "

# Row 24
$ws.Range("D24").Value = "13.34"
$ws.Range("E24").Value = "  +4.73%  "

# Row 25
$ws.Range("D25").Value = "3.21"
$ws.Range("E25").Value = "  +3.34%  "

# Row 26
$ws.Range("D26").Value = "32.77"
$ws.Range("E26").Value = "  +11.28%  "

# Row 27
$ws.Range("D27").Value = "8.81"
$ws.Range("E27").Value = "  +6.99%  "

# Row 28
$ws.Range("E28").Value = "  +0.32%  "

# Row 29
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").Value = "  -0.11%  "

# Row 30
$ws.Range("D30").Value = "2.73"
$ws.Range("E30").Value = "  +1.47%  "

# Row 31
$ws.Range("E31").Value = "  +1.91%  "

# Row 32
$ws.Range("E32").Value = "  +0.28%  "

# Row 33
$ws.Range("D33").Value = "43.67"
$ws.Range("E33").Value = "  +5.55%  "

# Row 34
$ws.Range("D34").Value = "11.85"

# Row 35
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("D36").Value = "0.0505"
$ws.Range("E36").Value = "  +5.58%  "

# Row 37
$ws.Range("D37").Value = "53.65"
$ws.Range("E37").Value = "  +3.84%  "

# Row 38
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  -0.16%  "

# Row 40
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -0.53%  "

# Row 41
$ws.Range("E41").Value = "  +6.42%  "

# Row 42
$ws.Range("E42").Value = "  +7.34%  "

# Row 43
$ws.Range("D43").Value = "142.31"
$ws.Range("E43").Value = "  +2.28%  "

# Row 44
$ws.Range("E44").Value = "  +0.41%  "

# Row 45
$ws.Range("E45").Value = "  +2.09%  "

# Row 46
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +9.01%  "

# Row 47
$ws.Range("D47").Value = "16.64"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("D48").Value = "21.88"
$ws.Range("E48").Value = "  +3.14%  "

# Row 49
$ws.Range("D49").Value = "2.121.32"
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("E50").Value = "  +16.01%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  +3.73%  "
